$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3c"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.227933
$ws.Range("H2").Value = 0.683799
$ws.Range("I2").Value = 0.005274728560320285
$ws.Range("J2").Value = 0.005274728560320286
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 12.06000424565433
$ws.Range("R2").Value = 108.540038210889
$ws.Range("S2").Value = 0.00219503367965828
$ws.Range("T2").Value = 0.00219503367965828

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3c"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.227933
$ws.Range("H3").Value = 0.683799
$ws.Range("I3").Value = 0.005274728560320285
$ws.Range("J3").Value = 0.005274728560320286
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 10.77212687609167
$ws.Range("R3").Value = 96.949141884825
$ws.Range("S3").Value = 0.001960627941162921
$ws.Range("T3").Value = 0.001960627941162921

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3c"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.227933
$ws.Range("H4").Value = 0.683799
$ws.Range("I4").Value = 0.005274728560320285
$ws.Range("J4").Value = 0.005274728560320286
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 6.148403173308668
$ws.Range("R4").Value = 55.335628559778
$ws.Range("S4").Value = 0.001119066939499085
$ws.Range("T4").Value = 0.001119066939499085

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema3c"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 42.071953
$ws.Range("H5").Value = 126.215859
$ws.Range("I5").Value = 0.9736112457500787
$ws.Range("J5").Value = 0.9736112457500787
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 2226.039809079728
$ws.Range("R5").Value = 20034.35828171755
$ws.Range("S5").Value = 0.4051600856567509
$ws.Range("T5").Value = 0.4051600856567509

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3c"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 42.071953
$ws.Range("H6").Value = 126.215859
$ws.Range("I6").Value = 0.9736112457500787
$ws.Range("J6").Value = 0.9736112457500787
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 1988.322952977258
$ws.Range("R6").Value = 17894.90657679532
$ws.Range("S6").Value = 0.3618933923174493
$ws.Range("T6").Value = 0.3618933923174493

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3c"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 42.071953
$ws.Range("H7").Value = 126.215859
$ws.Range("I7").Value = 0.9736112457500787
$ws.Range("J7").Value = 0.9736112457500787
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 1134.874411921455
$ws.Range("R7").Value = 10213.8697072931
$ws.Range("S7").Value = 0.2065577677758786
$ws.Range("T7").Value = 0.2065577677758786

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema3c"
$ws.Range("C8").Value = "Nrp1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.912385
$ws.Range("H8").Value = 2.737155
$ws.Range("I8").Value = 0.02111402568960099
$ws.Range("J8").Value = 0.021114025689601
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 48.27456741091166
$ws.Range("R8").Value = 434.471106698205
$ws.Range("S8").Value = 0.008786423220047204
$ws.Range("T8").Value = 0.008786423220047206

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema3c"
$ws.Range("C9").Value = "Nrp1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.912385
$ws.Range("H9").Value = 2.737155
$ws.Range("I9").Value = 0.02111402568960099
$ws.Range("J9").Value = 0.021114025689601
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 43.11936832245833
$ws.Range("R9").Value = 388.074314902125
$ws.Range("S9").Value = 0.007848128722466388
$ws.Range("T9").Value = 0.00784812872246639

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema3c"
$ws.Range("C10").Value = "Nrp1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.912385
$ws.Range("H10").Value = 2.737155
$ws.Range("I10").Value = 0.02111402568960099
$ws.Range("J10").Value = 0.021114025689601
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 24.61122711182334
$ws.Range("R10").Value = 221.50104400641
$ws.Range("S10").Value = 0.004479473747087401
$ws.Range("T10").Value = 0.004479473747087402

Write-Output "Done"
